# sd_covid_dataset.xlsx update: "Data from Apr 3"
#
# 1. Four new city columns (alpine, borregosprings, descando, valleycenter)
#    are inserted into the already-alphabetised city block (AV..BE before
#    the edit), shifting everything from their insertion point rightward.
# 2. tested (column B) gets a new value for Apr 2 (row 29).
# 3. A brand-new row (30) is appended for Apr 3, 2020.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the 4 new city columns, left to right -----------------------
# Inserting before AV shifts bonita..other (AV..BE) one column right.
$ws.Range("AV1").EntireColumn.Insert()
$ws.Range("AV1").Value = "alpine"

# Inserting before (then-)AX shifts jamul..other one column right; AW (bonita)
# is untouched.
$ws.Range("AX1").EntireColumn.Insert()
$ws.Range("AX1").Value = "borregosprings"

# Inserting before (then-)AY shifts lakeside..other one column right.
$ws.Range("AY1").EntireColumn.Insert()
$ws.Range("AY1").Value = "descando"

# Inserting before (then-)BH shifts other one column right (to BI).
$ws.Range("BH1").EntireColumn.Insert()
$ws.Range("BH1").Value = "valleycenter"

# --- 2. New "tested" figure for Apr 2 (row 29) ------------------------------
$ws.Range("B29").Value = 1882

# --- 3. New row 30 (Apr 3, 2020 = serial 43924) -----------------------------
$ws.Range("A30").Value = 43924
$ws.Range("C30").Value = 1112
$ws.Range("D30").Value = 10
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 12
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 198
$ws.Range("I30").Value = 11
$ws.Range("J30").Value = 238
$ws.Range("K30").Value = 25
$ws.Range("L30").Value = 192
$ws.Range("M30").Value = 29
$ws.Range("N30").Value = 177
$ws.Range("O30").Value = 35
$ws.Range("P30").Value = 136
$ws.Range("Q30").Value = 41
$ws.Range("R30").Value = 88
$ws.Range("S30").Value = 34
$ws.Range("T30").Value = 55
$ws.Range("U30").Value = 33
$ws.Range("V30").Value = 6
$ws.Range("W30").Value = 1
$ws.Range("X30").Value = 519
$ws.Range("Y30").Value = 585
$ws.Range("Z30").Value = 8
$ws.Range("AA30").Value = 211
$ws.Range("AB30").Value = 85
$ws.Range("AC30").Value = 17
$ws.Range("AD30").Value = 37
$ws.Range("AE30").Value = 75
$ws.Range("AF30").Value = 4
$ws.Range("AG30").Value = 7
$ws.Range("AH30").Value = 60
$ws.Range("AI30").Value = 27
$ws.Range("AJ30").Value = 26
$ws.Range("AK30").Value = 2
$ws.Range("AL30").Value = 15
$ws.Range("AM30").Value = 10
$ws.Range("AN30").Value = 18
$ws.Range("AO30").Value = 27
$ws.Range("AP30").Value = 13
$ws.Range("AQ30").Value = 614
$ws.Range("AR30").Value = 13
$ws.Range("AS30").Value = 12
$ws.Range("AT30").Value = 5
$ws.Range("AU30").Value = 15
$ws.Range("AV30").Value = 1
$ws.Range("AW30").Value = 7
$ws.Range("AX30").Value = 1
$ws.Range("AY30").Value = 1
$ws.Range("AZ30").Value = 5
$ws.Range("BA30").Value = 1
$ws.Range("BB30").Value = 8
$ws.Range("BC30").Value = 2
$ws.Range("BD30").Value = 6
$ws.Range("BE30").Value = 13
$ws.Range("BG30").Value = 23
$ws.Range("BH30").Value = 1
$ws.Range("BI30").Value = 46

# --- 4. Housekeeping that mirrors what Excel itself would update ----------
# The hidden AutoFilter-database defined name grows with the (now wider)
# header row.
foreach ($n in $wb.Names) {
  if ($n.Name -eq "Sheet1!_FilterDatabase") {
    $n.RefersTo = "=Sheet1!`$A`$1:`$BE`$17"
  }
}

# Selection ends up on the last cell touched (Y30), matching the saved view.
$ws.Range("Y30").Select()
